$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New text content (added in this order so the shared-string table
#     gets the same append order as the target file: "think this is
#     complete..." before "tom" before the Boss status text) ---
$ws.Range("D8").Value = "think this is complete now, need to actually test with another controller though"
$ws.Range("D9").Value = "think this is complete now, need to actually test with another controller though"
$ws.Range("B8").Value = "tom"
$ws.Range("B9").Value = "tom"

# Boss row: clear the old placeholder "completerer" and give it a status
$ws.Range("B12").ClearContents()
$ws.Range("D12").Value = "created class for boss, has health that player can decrease with sword swing, needs delay, is currently broken from players end, calls the lose hp for every frame of swing"

# --- Highlight the "done / in-progress" rows with a green fill ---
# (seed the green color on the Boss status cell first too, so the
#  subsequent ThemeColor call on the same cell re-uses this fill slot
#  instead of minting a throwaway intermediate fill)
$ws.Range("D12").Interior.Color = 5296274
$ws.Range("A2:D2").Interior.Color = 5296274
$ws.Range("A4:D9").Interior.Color = 5296274

# Boss status cell gets a plain white (Background 1 theme) fill instead
$ws.Range("D12").Interior.ThemeColor = 2

# --- Selection moves to B15 ---
$ws.Range("B15").Select()
